# The underlying data set gained one additional weekly observation.
# A new record was inserted at row 223 (pushing the existing rows 223:297
# down to 224:298). The new row duplicates the row immediately above it
# (row 222 - same market/product/quality/prices) but carries a new
# "Fecha" (date serial) value of 44524.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 223, shifting rows 223:297 down to 224:298.
$ws.Rows.Item(223).Insert()

# Populate the new row 223 with a copy of row 222's values.
$srcRange = $ws.Range("A222:R222")
$dstRange = $ws.Range("A223:R223")
$dstRange.Value2 = $srcRange.Value2

# The new row records a later date than its template row.
$ws.Range("D223").Value2 = 44524
